$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Direct value edits on existing rows (no shift) ---
    $ws.Range("F2").Value = 36
    $ws.Range("F3").Value = 195
    $ws.Range("F5").Value = 249
    $ws.Range("F6").Value = 38
    $ws.Range("F7").Value = 132

    # --- Insert a new row at position 9 (shifts old rows 9-35 down to 10-36) ---
    $ws.Rows.Item(9).Insert()

    # Copy column-A formatting (bold/border/centered) from the row below onto the new row
    $ws.Range("A10").Copy()
    $ws.Range("A9").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # --- Renumber the sequential index in column A for rows 9-36 (0-based: row9 -> 8, row36 -> 35) ---
    for ($r = 9; $r -le 36; $r++) {
        $ws.Range("A" + $r).Value = $r - 1
    }

    # --- Populate the newly inserted row 9 ---
    $ws.Range("A9").Value = 8
    $ws.Range("B9").Value = "'" + "2024-06-30"
    $ws.Range("C9").Value = "南昌·ChinastyleCOSPLAY  "
    $ws.Range("D9").Value = "真君路999号 南昌玛雅乐园"
    $ws.Range("E9").Value = "2024.06.30 09:30-07.02 17:30"
    $ws.Range("F9").Value = 4
    $ws.Range("G9").Value = 65
    $ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=87045"
    $ws.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202405/wajWy7ID1717149642528.jpeg"

    # --- Update want-to-go (F) / min-price (G) figures on the shifted rows 10-36 ---
    $ws.Range("F10").Value = 246
    $ws.Range("G10").Value = 55
    $ws.Range("F11").Value = 14
    $ws.Range("G11").Value = "不可售"
    $ws.Range("F12").Value = 36
    $ws.Range("G12").Value = 55
    $ws.Range("F13").Value = 23
    $ws.Range("G13").Value = 55
    $ws.Range("F14").Value = 82
    $ws.Range("G14").Value = 55
    $ws.Range("F15").Value = 394
    $ws.Range("G15").Value = 65
    $ws.Range("F16").Value = 41
    $ws.Range("G16").Value = 55
    $ws.Range("F17").Value = 472
    $ws.Range("G17").Value = 52.1
    $ws.Range("F18").Value = 397
    $ws.Range("G18").Value = 55
    $ws.Range("F19").Value = 134
    $ws.Range("G19").Value = 52.5
    $ws.Range("F20").Value = 62
    $ws.Range("G20").Value = 30
    $ws.Range("F21").Value = 31
    $ws.Range("G21").Value = 45
    $ws.Range("F22").Value = 38
    $ws.Range("G22").Value = 40
    $ws.Range("F23").Value = 1036
    $ws.Range("G23").Value = 19.9
    $ws.Range("F24").Value = 2786
    $ws.Range("G24").Value = 69
    $ws.Range("F25").Value = 22
    $ws.Range("G25").Value = 60
    $ws.Range("F26").Value = 55
    $ws.Range("G26").Value = 56
    $ws.Range("F27").Value = 530
    $ws.Range("G27").Value = 64
    $ws.Range("F28").Value = 33
    $ws.Range("G28").Value = 9.9
    $ws.Range("F29").Value = 971
    $ws.Range("G29").Value = 55
    $ws.Range("F30").Value = 565
    $ws.Range("G30").Value = "已售罄"
    $ws.Range("F31").Value = 451
    $ws.Range("G31").Value = 45
    $ws.Range("F32").Value = 260
    $ws.Range("G32").Value = 45
    $ws.Range("F33").Value = 387
    $ws.Range("G33").Value = 55
    $ws.Range("F34").Value = 450
    $ws.Range("G34").Value = 45
    $ws.Range("F35").Value = 594
    $ws.Range("G35").Value = 45
    $ws.Range("F36").Value = 422
    $ws.Range("G36").Value = 45
}

Write-Output "done"
